# Extend European routes with a new Faro -> Brussels entry.
# A new row is inserted at row 22 (pushing the existing IndianWings rows
# down to 23-25), matching the EuropeanWings block already present in
# rows 17-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 22, shifting rows 22+ down.
$ws.Rows("22:22").Insert()

# Fill in the new route. Order matches how the shared-string table grew
# in the source edit (ICAO arrival, departure name, arrival name, ICAO
# departure, then airline).
$ws.Range("E22").Value = "EBBR"
$ws.Range("B22").Value = "Faro-Portugal"
$ws.Range("D22").Value = "Brussels-National"
$ws.Range("C22").Value = "LPFR"
$ws.Range("A22").Value = "EuropeanWings"

# Leave the selection where the author left it when saving.
$ws.Range("B30").Select() | Out-Null
